$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$ws1.Range("D5").Value = "2016-03-31 06:36:58"
$ws2.Range("E5").Value = "2016-03-31 06:36:48"
$ws3.Range("E5").Value = "2016-03-31 06:36:58"
